$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 59, shifting the existing data (rows 59-94) down to rows 61-96
$ws.Rows("59:60").Insert()

# New row 59: Primera quality entry dated 44634
$ws.Range("A59").Value = 1
$ws.Range("B59").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C59").Value = "Arica y Parinacota"
$ws.Range("D59").Value = 44634
$ws.Range("E59").Value = 15
$ws.Range("F59").Value = 100112008
$ws.Range("G59").Value = "Coliflor"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 700
$ws.Range("K59").Value = 1100
$ws.Range("L59").Value = 1200
$ws.Range("M59").Value = 1150
$ws.Range("N59").Value = "$/unidad"
$ws.Range("O59").Value = "Región de Arica y Parinacota"
$ws.Range("P59").Value = 1150
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"

# New row 60: Segunda quality entry dated 44634
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44634
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = 100112008
$ws.Range("G60").Value = "Coliflor"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Segunda"
$ws.Range("J60").Value = 800
$ws.Range("K60").Value = 800
$ws.Range("L60").Value = 900
$ws.Range("M60").Value = 850
$ws.Range("N60").Value = "$/unidad"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 850
$ws.Range("Q60").Value = 1
$ws.Range("R60").Value = "Hortaliza"
